# Apply the "Optuna Attempt (go back with original)" revisions.
# Sheet "Forecast Comparison" (rows 2-17): update MyForecast (D),
# Inventory Coverage (H), Stockout Risk (I) and Seasonality Index (L).
# Sheet "Summary": update the rolled-up forecast totals (B9, B10, B11,
# B12, B14) -- these cells hold numeric-looking text, so they are written
# with a leading apostrophe to keep them stored as text, matching the
# original file.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------------

# Row 2 (W8)
$wsForecast.Range("D2").Value = 169
$wsForecast.Range("H2").Value = 2.5
$wsForecast.Range("L2").Value = 1.09

# Row 3 (W9)
$wsForecast.Range("D3").Value = 185
$wsForecast.Range("H3").Value = 1.37
$wsForecast.Range("L3").Value = 1.05

# Row 4 (W10)
$wsForecast.Range("D4").Value = 154
$wsForecast.Range("H4").Value = 0.44
$wsForecast.Range("I4").Value = "High"
$wsForecast.Range("L4").Value = 0.99

# Row 5 (W11)
$wsForecast.Range("L5").Value = 0.88

# Row 6 (W12)
$wsForecast.Range("L6").Value = 1.15

# Row 7 (W13)
$wsForecast.Range("L7").Value = 0.91

# Row 8 (W14)
$wsForecast.Range("D8").Value = 80
$wsForecast.Range("L8").Value = 0.83

# Row 9 (W15)
$wsForecast.Range("D9").Value = 79
$wsForecast.Range("L9").Value = 0.95

# Row 10 (W16)
$wsForecast.Range("D10").Value = 77
$wsForecast.Range("L10").Value = 0.8100000000000001

# Row 11 (W17)
$wsForecast.Range("D11").Value = 75
$wsForecast.Range("L11").Value = 0.98

# Row 12 (W18)
$wsForecast.Range("D12").Value = 74
$wsForecast.Range("L12").Value = 1.05

# Row 13 (W19)
$wsForecast.Range("D13").Value = 75
$wsForecast.Range("L13").Value = 0.91

# Row 14 (W20)
$wsForecast.Range("D14").Value = 89
$wsForecast.Range("L14").Value = 0.86

# Row 15 (W21)
$wsForecast.Range("D15").Value = 69
$wsForecast.Range("L15").Value = 1.13

# Row 16 (W22)
$wsForecast.Range("D16").Value = 70
$wsForecast.Range("L16").Value = 0.91

# Row 17 (W23)
$wsForecast.Range("D17").Value = 67
$wsForecast.Range("L17").Value = 1.15

# --- Summary sheet --------------------------------------------------------

$wsSummary.Range("B9").Value  = "'1733"
$wsSummary.Range("B10").Value = "'1134"
$wsSummary.Range("B11").Value = "'667"
$wsSummary.Range("B12").Value = "'186"
$wsSummary.Range("B14").Value = "'68"
